# Update "想去人数" (interested-count) values that changed between crawls.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 2041
$ws1.Range("F19").Value = 115
$ws1.Range("F21").Value = 238
$ws1.Range("F24").Value = 415
$ws1.Range("F25").Value = 230
$ws1.Range("F27").Value = 341

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 2041
$ws4.Range("F26").Value = 115
$ws4.Range("F28").Value = 238
$ws4.Range("F31").Value = 415
$ws4.Range("F34").Value = 230
$ws4.Range("F36").Value = 341
